$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Job to Run")
$ws2 = $wb.Worksheets.Item("All")

# Remove the "Exp Folder Name" and "CTRL Folder Name" columns from both
# sheets (script1 now drops CTRL/EXP .mzML files into a single shared
# folder per job, so the separate folder-name columns are no longer
# needed).  Delete the later column first so the earlier column index
# stays valid.
$ws1.Range("D1:D2").EntireColumn.Delete() | Out-Null
$ws1.Range("B1:B2").EntireColumn.Delete() | Out-Null

$ws2.Range("D1:D2").EntireColumn.Delete() | Out-Null
$ws2.Range("B1:B2").EntireColumn.Delete() | Out-Null

# On the "Job to Run" sheet, use the full job/experiment name (as already
# used on the "All" sheet) instead of the old short alias.
$ws2.Range("A2").Copy($ws1.Range("A2")) | Out-Null

$ws2.Range("A2").Select() | Out-Null
$ws1.Range("C10").Select() | Out-Null
